$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.295.06'
$ws.Range("E2").Value = '  +8.53%  '
$ws.Range("D3").Value = '1.596.87'
$ws.Range("E3").Value = '  +8.32%  '
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'0.9937"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.04%  '
$ws.Range("D6").Value = "'302.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.73%  '
$ws.Range("D7").Value = "'0.3635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").Value = "'0.3378"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.08%  '
$ws.Range("D9").Value = "'41.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.99%  '
$ws.Range("D10").Value = "'1.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").Value = "'0.06998"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.10%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = "'19.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.61%  '
$ws.Range("D14").Value = "'5.870"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.17%  '
$ws.Range("D15").Value = "'6.583"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.26%  '
$ws.Range("D16").Value = "'0.9937"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("D17").Value = "'0.00001066"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.76%  '
$ws.Range("D18").Value = '1.595.77'
$ws.Range("E18").Value = '  +8.21%  '
$ws.Range("D19").Value = "'0.06604"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.44%  '
$ws.Range("D20").Value = "'76.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.24%  '
$ws.Range("D21").Value = "'5.956"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.10%  '
$ws.Range("D22").Value = "'15.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.00%  '
$ws.Range("D23").Value = "'11.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.72%  '
$ws.Range("D24").Value = '22.365.98'
$ws.Range("E24").Value = '  +8.75%  '
$ws.Range("D25").Value = "'2.393"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.31%  '
$ws.Range("D26").Value = "'2.526"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +18.23%  '
$ws.Range("D27").Value = "'148.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.83%  '
$ws.Range("D28").Value = "'19.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.55%  '
$ws.Range("D29").Value = '1.767.07'
$ws.Range("E29").Value = '  +7.97%  '
$ws.Range("D30").Value = "'122.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.51%  '
$ws.Range("D31").Value = "'4.055"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.00%  '
$ws.Range("D32").Value = "'6.002"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +20.34%  '
$ws.Range("D33").Value = "'0.9294"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.73%  '
$ws.Range("D34").Value = "'1.668"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.15%  '
$ws.Range("D35").Value = "'0.08188"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.48%  '
$ws.Range("D36").Value = "'11.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.17%  '
$ws.Range("D37").Value = "'5.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.68%  '
$ws.Range("D38").Value = "'1.244"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("D39").Value = "'8.403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.14%  '
$ws.Range("D40").Value = "'0.06009"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.57%  '
$ws.Range("D41").Value = "'0.02190"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.51%  '
$ws.Range("D42").Value = "'0.1996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.15%  '
$ws.Range("D43").Value = "'0.9934"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.90%  '
$ws.Range("D44").Value = "'0.5825"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.18%  '
$ws.Range("D45").Value = "'3.802"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.94%  '
$ws.Range("D46").Value = "'13.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.83%  '
$ws.Range("D47").Value = "'0.5615"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.87%  '
$ws.Range("D48").Value = "'125.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.03%  '
$ws.Range("D49").Value = "'1.949"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.32%  '
$ws.Range("D50").Value = "'0.06766"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.43%  '
$ws.Range("D51").Value = "'72.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.01%  '